# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-06 09:14:35
#
# Applies the CVS session-analysis refresh:
#   - reordered / refreshed "Recorded By" (column G) email lists
#   - updated session + per-group statistics (columns L, P, Q, S)
#   - two sessions (rows 10 & 113) flipped from "Pending" to
#     "Not Recorded", including the row shading that goes with it
#
# Notes on technique:
#   * Percent-looking text such as "33.4%" gets silently reinterpreted
#     as a number by a plain `.Value =` assignment (Excel's "looks like
#     a percentage" autodetection), which would also change the cell's
#     number format / style. To keep these as literal text (matching
#     the original inline-string cells) we write a tiny text formula
#     and immediately collapse it back to a static value in place via
#     Copy / PasteSpecial(xlPasteValues) - this never touches
#     NumberFormat, so the existing style id is preserved exactly.
#   * Rows 10 and 113 need to switch from the "Pending" fill to the
#     "Not Recorded" fill. Rather than poking Interior.Color (which
#     would fabricate a brand-new style), we copy the *formatting only*
#     from row 7 - an existing "Not Recorded" row that already uses the
#     desired style - via Copy / PasteSpecial(xlPasteFormats).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    # Forces a literal text value into $range without Excel's numeric/
    # percent auto-detection touching the cell's format/style.
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = 0
}

# ---------------------------------------------------------------------
# Recorded-By (column G) email list reorderings / updates
# ---------------------------------------------------------------------
$ws.Range("G2").Value  = "nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G24").Value = "nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

$ws.Range("G18").Value = "aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G40").Value = "aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"

$ws.Range("G19").Value  = "naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G41").Value  = "naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G150").Value = "naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G172").Value = "naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

$ws.Range("G32").Value = "Madeha.Saeed@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("H32").Value = "55/217"

$ws.Range("G46").Value = "nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("G68").Value = "nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"

$ws.Range("G52").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G74").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

$emails_G54 = "yassmina.fattoh@med.asu.edu.eg, maimustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg"
$ws.Range("G54").Value  = $emails_G54
$ws.Range("G76").Value  = $emails_G54
$ws.Range("G98").Value  = $emails_G54
$ws.Range("G120").Value = $emails_G54
$ws.Range("H120").Value = "165/224"

$ws.Range("G62").Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G84").Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"

$emails_G96 = "Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G96").Value  = $emails_G96
$ws.Range("G118").Value = $emails_G96

$ws.Range("G106").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G128").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

$ws.Range("G107").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G129").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

$ws.Range("G134").Value = "Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

$emails_G142 = "yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg"
$ws.Range("G142").Value = $emails_G142
$ws.Range("G164").Value = $emails_G142

$ws.Range("G156").Value = "alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Class statistics (column L) on the summary block
# ---------------------------------------------------------------------
$ws.Range("L7").Value = 13
$ws.Range("L8").Value = 122
Set-TextValue $ws.Range("L10") "33.4%"

# ---------------------------------------------------------------------
# Per-group statistics block (columns P, Q, S)
# ---------------------------------------------------------------------
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 15
Set-TextValue $ws.Range("S16") "34.8%"

$ws.Range("P20").Value = 3
$ws.Range("Q20").Value = 14
Set-TextValue $ws.Range("S20") "42.2%"

# ---------------------------------------------------------------------
# Rows 10 and 113: sessions flipped from "Pending" to "Not Recorded".
# Copy the shading (fill) from an existing "Not Recorded" row (row 7)
# so the style index is reused rather than creating a brand-new one,
# then update the status text.
# ---------------------------------------------------------------------
$ws.Range("A7:I7").Copy()
$ws.Range("A10:I10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("I10").Value = "Not Recorded"

$ws.Range("A7:I7").Copy()
$ws.Range("A113:I113").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("I113").Value = "Not Recorded"
